# Automatische test-sync: 2025-06-19 10:00:10
#
# 1. Append two new e-mail log rows (4 and 5) on the "Logs" sheet.
# 2. Extend the conditional-formatting ranges on "Logs" (D and G columns)
#    to cover the two new rows.
# 3. Insert a new "Overig" category row on the "Dashboard" sheet (becomes
#    row 2), pushing "Afmelding" and "Klacht" down one row each.
# 4. Update the bar chart's category/value series references on the
#    "Dashboard" sheet so they include the extra row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: add the two new rows
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Sollicitatie marketingfunctie"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D4").Value = "Overig"
$logs.Range("F4").Value = "2025-06-19 09:58:11"
$logs.Range("G4").Value = "Nee"

$logs.Range("A5").Value = "Vragen over samenwerking"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D5").Value = "Overig"
$logs.Range("F5").Value = "2025-06-19 09:58:11"
$logs.Range("G5").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting ranges to include the
#    new rows (D2:D3 -> D2:D5, G2:G3 -> G2:G5)
# ---------------------------------------------------------------------
$dConditions = $logs.Range("D2:D3").FormatConditions
$dConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))

$gConditions = $logs.Range("G2:G3").FormatConditions
$gConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: add an "Overig" row as the new row 2, shifting the
#    existing "Afmelding"/"Klacht" rows down (manual shift keeps the
#    original, unstyled row formatting instead of inheriting the header
#    row's style the way Rows.Insert() would).
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = $dash.Range("A3").Value2
$dash.Range("B4").Value = $dash.Range("B3").Value2

$dash.Range("A3").Value = $dash.Range("A2").Value2
$dash.Range("B3").Value = $dash.Range("B2").Value2

$dash.Range("A2").Value = "Overig"
$dash.Range("B2").Value = 2

# ---------------------------------------------------------------------
# 4. Dashboard sheet: extend the chart series references so the chart
#    picks up the new "Overig" category row too.
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$4"
$series.Values = "='Dashboard'!`$B`$2:`$B`$4"
